$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Copy formatting from row 11 (A11:C11) into the new row 12 before filling values
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A12").Value = 43907
$ws.Range("B12").Value = 2.25
$ws.Range("C12").Value = "Préparation de l'environement de travail à la maison (sans WAMP) avec l'aide de Dorian"

$ws.Range("C13").Select()
